$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.902.54"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "2.686.04"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("E11").Value = "  -3.19%  "
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000205"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.23%  "
$ws.Range("D15").Value = "3.168.25"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").Value = "65.688.14"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "2.679.90"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "353.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.48%  "
$ws.Range("E24").Value = "  +9.66%  "
$ws.Range("E25").Value = "  +2.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.171"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.03%  "
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "531.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.56%  "
$ws.Range("E33").Value = "  -3.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.429"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "166.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0627"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0265"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.652"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0989"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.07%  "
